# The glossary table's "Definição" (definition) column had its font size
# reduced from 14pt (w:sz/w:szCs = 28) to 12pt (w:sz/w:szCs = 24) for the
# "Gateway" and "CTC" entries. Apply the same change via the Word object
# model by resizing the font of every paragraph in those two table cells.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Column 2 is the "Definição" column; rows 2 and 3 are "Gateway" and "CTC".
$targetRows = 2, 3

foreach ($rowIndex in $targetRows) {
    $cell = $table.Cell($rowIndex, 2)
    $cellRange = $cell.Range
    for ($p = 1; $p -le $cellRange.Paragraphs.Count; $p++) {
        $para = $cellRange.Paragraphs.Item($p)
        $para.Range.Font.Size = 12
        $para.Range.Font.SizeBi = 12
    }
}
